# Katalog guncellendi - Per 27.11.2025 10:34:41,30
# Adds 4 new "Sweat" category products (MASTIF 3050 3 Ip Sweatshirt) to the
# product catalog sheet, right after the existing last row (105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aciklama = "S-M-L-XL-2XL-3XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."

$rows = @(
    @("MASTİF 3050  3 İP SWEATSHİRT SİYAH",  "280 TL", "Sweat", "MASTİF3050SİY.jpg",   $aciklama, "var"),
    @("MASTİF 3050  3 İP SWEATSHİRT MAVİ",   "280 TL", "Sweat", "MASTİF3050MAVİ.jpg",  $aciklama, "var"),
    @("MASTİF 3050  3 İP SWEATSHİRT LACİVERT","280 TL", "Sweat", "MASTİF3050LACİ.jpg",  $aciklama, "var"),
    @("MASTİF 3050  3 İP SWEATSHİRT SU YEŞİLİ","280 TL", "Sweat", "MASTİF3050SUYEŞ.jpg", $aciklama, "var")
)

$startRow = 106
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    # Column entry order mirrors how the catalog was actually typed in:
    # price, category, image file, product name, description, stock.
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

$ws.Range("E102").Select()
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 1
